$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.732.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.322.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.37%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.19"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.02%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.94%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.957"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.320.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.38"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.12"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.633.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.938.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.03"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.321.03"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.23"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.42"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "490.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.438"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -9.88%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "92.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.89"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.499.36"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.10"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.61"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.01"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "28.20"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.526"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "562.54"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.32%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.865"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "23.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.67"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.39"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.98"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.79"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.11%  "
